$d = $word.ActiveDocument

# Locate the paragraph "I am almost done with this :-" (the last piece of
# real content before the document's trailing empty paragraph) and add a
# brand-new paragraph right after it that reads "Creared the further
# concept". "Creared" is kept as its own run, mirroring the misspelling
# Word's proofer flags separately in the source diff.
$anchor = $d.Content
$found = $anchor.Find.Execute("I am almost done with this :-", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()

$newPara1 = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$newPara1.Range.InsertAfter("Creared")

# Split "Creared" from " the further concept" into two distinct runs by
# temporarily making them separate paragraphs, then merging the
# paragraph mark back out. This keeps each piece of text in its own
# <w:r> run without leaving stray direct-formatting behind.
$newPara1.Range.InsertParagraphAfter()
$newPara2 = $newPara1.Next()
$newPara2.Range.InsertAfter(" the further concept")

$p1End = $newPara1.Range.End
$joinMark = $d.Range($p1End - 1, $p1End)
$joinMark.Delete()
